$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# and the original text formatting (e.g. trailing zeros) would be lost.
$forceTextCells = @("D5", "D8", "D9", "D17", "D18", "D19", "D26", "D27", "D33", "D36", "D38", "D40", "D44", "D47", "D49")
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "28.530.07"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.565.70"
$ws.Range("E3").Value = "  -2.01%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "211.84"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "46.16"
$ws.Range("E8").Value = "  +4.73%  "
$ws.Range("D9").Value = "24.05"
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  -1.74%  "
$ws.Range("E11").Value = "  -1.57%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "1.789.59"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "1.567.45"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "28.518.70"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "3.68"
$ws.Range("E17").Value = "  -3.40%  "
$ws.Range("D18").Value = "62.21"
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").Value = "229.03"
$ws.Range("E19").Value = "  -1.59%  "
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -5.94%  "
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("E25").Value = "  +6.51%  "
$ws.Range("D26").Value = "150.52"
$ws.Range("E26").Value = "  -1.54%  "
$ws.Range("D27").Value = "14.99"
$ws.Range("E27").Value = "  -2.02%  "
$ws.Range("E28").Value = "  -2.82%  "
$ws.Range("E29").Value = "  -4.14%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("E32").Value = "  -4.09%  "
$ws.Range("D33").Value = "3.21"
$ws.Range("E33").Value = "  -1.33%  "
$ws.Range("D35").Value = "1.392.04"
$ws.Range("E35").Value = "  -2.13%  "
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").Value = "2.36"
$ws.Range("E38").Value = "  +0.84%  "
$ws.Range("E39").Value = "  +1.92%  "
$ws.Range("D40").Value = "0.0166"
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  -1.76%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("D44").Value = "0.789"
$ws.Range("E44").Value = "  -4.07%  "
$ws.Range("E45").Value = "  -4.31%  "
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "62.80"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").Value = "1.702.18"
$ws.Range("E48").Value = "  -1.98%  "
$ws.Range("D49").Value = "86.08"
$ws.Range("E49").Value = "  -1.67%  "
$ws.Range("E50").Value = "  -4.37%  "
$ws.Range("E51").Value = "  -0.43%  "
